# penambahan data lokasi pegawai
# - remove the erroneous "Nganjuk" row (id 2) from the city list
# - append a new city "Cimahi" (id 23) at the end of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 currently holds ID=2 / Nama Kota="Nganjuk" - delete it, shifting
# the rows below it up by one.
$ws.Rows.Item(3).Delete()

# Append the new city as the new last row of the table (row 8 after the
# delete above).
$ws.Cells.Item(8, 1).Value = 23
$ws.Cells.Item(8, 2).Value = "Cimahi"
